$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 439.33334
$ws.Range("I28").Value = 459.875
$ws.Range("K28").Value = 459.875
$ws.Range("M28").Value = 25.125
$ws.Range("H32").Value = 9434.5
$ws.Range("I32").Value = 10000
$ws.Range("K32").Value = 10000
$ws.Range("M32").Value = -9674
$ws.Range("H62").Value = 6166
$ws.Range("I62").Value = 5999
$ws.Range("K62").Value = 5999
$ws.Range("M62").Value = -5375
$ws.Range("H65").Value = 6166
$ws.Range("I65").Value = 5999
$ws.Range("K65").Value = 29995
$ws.Range("M65").Value = -26875
$ws.Range("H69").Value = 1666.6666
$ws.Range("I69").Value = 1000
$ws.Range("K69").Value = 3000
$ws.Range("M69").Value = -2126
$ws.Range("H72").Value = 1666.6666
$ws.Range("I72").Value = 1000
$ws.Range("K72").Value = 9000
$ws.Range("M72").Value = -4632
$ws.Range("H76").Value = 4710
$ws.Range("I76").Value = 4710
$ws.Range("K76").Value = 4710
$ws.Range("M76").Value = -4395
$ws.Range("H79").Value = 4710
$ws.Range("I79").Value = 4710
$ws.Range("K79").Value = 4710
$ws.Range("M79").Value = -3618
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H100").Value = 11918.714
$ws.Range("J100").Value = 27500
$ws.Range("L100").Value = 27500
$ws.Range("N100").Value = -28582
$ws.Range("H111").Value = 3993.75
$ws.Range("I111").Value = 4660
$ws.Range("J111").Value = 1995
$ws.Range("K111").Value = 13980
$ws.Range("L111").Value = 5985
$ws.Range("M111").Value = -10913
$ws.Range("N111").Value = -12119
$ws.Range("H135").Value = 833.5714
$ws.Range("I135").Value = 733.25
$ws.Range("J135").Value = 1435.5
$ws.Range("K135").Value = 6599.25
$ws.Range("L135").Value = 12919.5
$ws.Range("M135").Value = -4064.25
$ws.Range("N135").Value = -17989.5
$ws.Range("H138").Value = 7442.3936
$ws.Range("J138").Value = 8016.2544
$ws.Range("L138").Value = 24048.7632
$ws.Range("N138").Value = -34328.7632

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4257.4
$ws.Range("I61").Value = 4916.909
$ws.Range("J61").Value = 2443.75
$ws.Range("K61").Value = 4916.909
$ws.Range("L61").Value = 2443.75
$ws.Range("M61").Value = -4704.909
$ws.Range("N61").Value = -2867.75
$ws.Range("H76").Value = 32500
$ws.Range("J76").Value = 32500
$ws.Range("L76").Value = 32500
$ws.Range("N76").Value = -33176
$ws.Range("H79").Value = 32500
$ws.Range("J79").Value = 32500
$ws.Range("L79").Value = 32500
$ws.Range("N79").Value = -34840
$ws.Range("H110").Value = 2118.75
$ws.Range("I110").Value = 1992
$ws.Range("J110").Value = 2499
$ws.Range("K110").Value = 1992
$ws.Range("L110").Value = 2499
$ws.Range("M110").Value = 53
$ws.Range("N110").Value = -6589
$ws.Range("H122").Value = 7338.933
$ws.Range("I122").Value = 7758.077
$ws.Range("J122").Value = 4614.5
$ws.Range("K122").Value = 23274.231
$ws.Range("L122").Value = 13843.5
$ws.Range("M122").Value = -20824.231
$ws.Range("N122").Value = -18743.5
$ws.Range("H132").Value = 3899
$ws.Range("I132").Value = 2799
$ws.Range("K132").Value = 8397
$ws.Range("M132").Value = -5867
$ws.Range("H136").Value = 4257.4
$ws.Range("I136").Value = 4916.909
$ws.Range("J136").Value = 2443.75
$ws.Range("K136").Value = 14750.727
$ws.Range("L136").Value = 7331.25
$ws.Range("M136").Value = -12200.727
$ws.Range("N136").Value = -12431.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3216.3333
$ws.Range("I20").Value = 3177.5
$ws.Range("J20").Value = 3294
$ws.Range("K20").Value = 3177.5
$ws.Range("L20").Value = 3294
$ws.Range("M20").Value = -2930.5
$ws.Range("N20").Value = -3788
$ws.Range("H105").Value = 1357.2858
$ws.Range("I105").Value = 1384.8462
$ws.Range("K105").Value = 1384.8462
$ws.Range("M105").Value = 362.1538
$ws.Range("H107").Value = 4968.125
$ws.Range("I107").Value = 4991.6
$ws.Range("K107").Value = 4991.6
$ws.Range("M107").Value = -3071.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 260
$ws.Range("J7").Value = 260
$ws.Range("L7").Value = 260
$ws.Range("N7").Value = -486
$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1700
$ws.Range("H94").Value = 2411.4285
$ws.Range("J94").Value = 2363.5
$ws.Range("L94").Value = 2363.5
$ws.Range("N94").Value = -3265.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 85
$ws.Range("I2").Value = 47.5
$ws.Range("J2").Value = 122.5
$ws.Range("K2").Value = 285
$ws.Range("L2").Value = 735
$ws.Range("M2").Value = -172
$ws.Range("N2").Value = -961
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H102").Value = 5000
$ws.Range("J102").Value = 5000
$ws.Range("L102").Value = 15000
$ws.Range("N102").Value = -19868
$ws.Range("H104").Value = 4845
$ws.Range("J104").Value = 4845
$ws.Range("L104").Value = 14535
$ws.Range("N104").Value = -19777
$ws.Range("H129").Value = 2068.4285
$ws.Range("I129").Value = 989.5
$ws.Range("K129").Value = 2968.5
$ws.Range("M129").Value = 2031.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H126").Value = 5744.25
$ws.Range("I126").Value = 3193.6
$ws.Range("J126").Value = 9995.333000000001
$ws.Range("K126").Value = 9580.799999999999
$ws.Range("L126").Value = 29985.999
$ws.Range("M126").Value = -7110.799999999999
$ws.Range("N126").Value = -34925.999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2856.6155
$ws.Range("I16").Value = 2557
$ws.Range("K16").Value = 2557
$ws.Range("M16").Value = -2387
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H100").Value = 1298.5
$ws.Range("I100").Value = 1298.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1298.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -757.5
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 3235.5
$ws.Range("I122").Value = 3246.3333
$ws.Range("K122").Value = 9738.999899999999
$ws.Range("M122").Value = -7288.999899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1053.4286
$ws.Range("I107").Value = 1080
$ws.Range("J107").Value = 1033.5
$ws.Range("K107").Value = 3240
$ws.Range("L107").Value = 3100.5
$ws.Range("M107").Value = -1320
$ws.Range("N107").Value = -6940.5
$ws.Range("H122").Value = 3959.2
$ws.Range("I122").Value = 3933.6667
$ws.Range("J122").Value = 3997.5
$ws.Range("K122").Value = 11801.0001
$ws.Range("L122").Value = 11992.5
$ws.Range("M122").Value = -9351.000100000001
$ws.Range("N122").Value = -16892.5
